$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "BEFORE insert Q45 style via NumberFormat:" $ws.Range("Q45").NumberFormat
$ws.Rows("45").Insert()
Write-Host "AFTER insert Q46 NumberFormat:" $ws.Range("Q46").NumberFormat
$ws.Range("A46:Q46").Copy()
$ws.Range("A45:Q45").PasteSpecial(-4122)
Write-Host "AFTER paste Q45 NumberFormat:" $ws.Range("Q45").NumberFormat
Write-Host "AFTER paste N45 NumberFormat:" $ws.Range("N45").NumberFormat
